# "Added padding for empty cell"
#
# Sheet1:
#   - Insert a new column before E, with header "Surprise" (old E..G shift to F..H).
#   - Insert a new blank padding row before row 3 (old row 3 becomes row 4).
#   - Type some "padding"/filler values into column I, rows 2 and 4.
#   - Hyperlinks on the Email column need to follow the row shift (C2 stays,
#     the old C3 hyperlink now belongs to C4).
#
# Test sheet: untouched (its lone cell keeps referencing the same text,
# only its underlying shared-string id moves because new strings were
# introduced on Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- insert the new "Surprise" column before the old column E ---
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("E1").Value = "Surprise"
$ws.Range("E1").ColumnWidth = 11

# --- insert a new blank padding row before the old row 3 ---
$ws.Range("A3").EntireRow.Insert()

# --- padding values typed into the new column I ---
$ws.Range("I2").Value = "dsdsqd"
$ws.Range("I4").Value = "qdqsdq"

# --- stash the original hyperlink-cell look (blue font, no underline) ---
$ws.Range("C2").Copy()
$ws.Range("Z100").PasteSpecial(-4122)

# --- rebuild the hyperlinks so the Lefebvre one now targets C4 ---
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:john.smith@mail.com", "", "", "john.smith@mail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:jean.lefebvre@mail.com", "", "", "jean.lefebvre@mail.com")

# --- Hyperlinks.Add overwrote the cell style with the built-in "Hyperlink"
#     style; restore the workbook's own custom hyperlink look ---
$ws.Range("Z100").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

# --- match the author's final selection ---
$ws.Range("C11").Select() | Out-Null
